# admission.xlsx — project structure rename + rule-engine font tidy-up
#
# 1) The Drools import rows (B2, B3) move to the new package layout:
#      com.admission.drools.api.Student              -> com.admission.drools.api.model.Student
#      com.admission.drools.api.StudentIdGenerator    -> com.admission.drools.api.utilities.StudentIdGenerator
#
# 2) The filler/spacer cells to the right of the rule table (columns F:H on
#    rows 1-10, and E:H on rows 11-12) use a 12pt font that never got a face
#    or color assigned. Give it an explicit Arial face + automatic (theme)
#    text color, matching the rest of the sheet's fonts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the Import statements -------------------------------------
$ws.Range("B2").Value = "com.admission.drools.api.model.Student"
$ws.Range("B3").Value = "com.admission.drools.api.utilities.StudentIdGenerator"

# --- 2) Give the unnamed 12pt font a face + color -------------------------
# NOTE: the headless Font setter reconstructs a font record from whatever it
# can currently read back, and the bare "sz 12" font (no name/color set)
# reads its size back incorrectly unless Size is (re)asserted alongside the
# new properties. Re-assert Size=12 so the rebuilt font keeps its 12pt size
# instead of silently collapsing to the default 10pt.
$fillerRanges = @("F1:H10", "E11:H12")
foreach ($addr in $fillerRanges) {
    $rng = $ws.Range($addr)
    $rng.Font.Size = 12
    $rng.Font.Name = "Arial"
    $rng.Font.ColorIndex = -4105   # xlColorIndexAutomatic -> <color theme="1"/>
}
